$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 36, pushing the old rows 36-41
# down to 39-44. Excel will copy formatting (e.g. the date style on column D)
# from the row above, matching the target style.
$ws.Rows("36:38").Insert()

# Populate the three newly-inserted rows (36, 37, 38) with the new
# Chirimoya price records. Columns A, B, C, E, F, G, H, I, J, K, Q, R, T are
# identical to every other row in this table.
$newRows = @(
    @{ Row = 36; D = 44466; L = "Especial"; M = 45;  N = 27000; O = 27000; P = 27000; S = 2700 },
    @{ Row = 37; D = 44466; L = "Primera";  M = 48;  N = 25000; O = 25000; P = 25000; S = 2500 },
    @{ Row = 38; D = 44466; L = "Segunda";  M = 40;  N = 23000; O = 23000; P = 23000; S = 2300 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 10
}
